# Generate Report for Handback
# Updates timestamp values recorded on the "Overview", "zh-cn" and "de-de"
# sheets to reflect the latest handoff/handback generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-17 21:06:17"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-17 21:06:12"
$wsZhCn.Range("K2").Value = "2016-08-17 21:06:30"

# de-de sheet: same two columns for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-17 21:06:17"
$wsDeDe.Range("K2").Value = "2016-08-17 21:06:39"
